$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1577
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1730.8572
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 5192.571599999999
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -5528.571599999999

$ws.Range("H74").Value = 270833
$ws.Range("I74").Value = 270833
$ws.Range("K74").Value = 270833
$ws.Range("M74").Value = -269897

$ws.Range("H77").Value = 270833
$ws.Range("I77").Value = 270833
$ws.Range("K77").Value = 1354165
$ws.Range("M77").Value = -1349485

$ws.Range("H86").Value = 8497.583000000001
$ws.Range("I86").Value = 8285.429
$ws.Range("K86").Value = 8285.429
$ws.Range("M86").Value = -7162.429

$ws.Range("H89").Value = 8497.583000000001
$ws.Range("I89").Value = 8285.429
$ws.Range("K89").Value = 41427.145
$ws.Range("M89").Value = -35811.145

$ws.Range("H116").Value = 20783.812
$ws.Range("I116").Value = 5266
$ws.Range("K116").Value = 5266
$ws.Range("M116").Value = -1824

$ws.Range("H132").Value = 5299.3584
$ws.Range("I132").Value = 5274.3125
$ws.Range("K132").Value = 15822.9375
$ws.Range("M132").Value = -13292.9375

$ws.Range("H137").Value = 15389000
$ws.Range("I137").Value = 35715920
$ws.Range("K137").Value = 107147760
$ws.Range("M137").Value = -107145210

$ws.Range("H138").Value = 2829.0908
$ws.Range("I138").Value = 1946.2778
$ws.Range("J138").Value = 3440.2693
$ws.Range("K138").Value = 5838.8334
$ws.Range("L138").Value = 10320.8079
$ws.Range("M138").Value = -698.8334000000004
$ws.Range("N138").Value = -20600.8079

$ws.Range("H141").Value = 7537.647
$ws.Range("I141").Value = 2414
$ws.Range("K141").Value = 7242
$ws.Range("M141").Value = -2062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 10244.833
$ws.Range("I31").Value = 10244.833
$ws.Range("K31").Value = 10244.833
$ws.Range("M31").Value = -9950.833000000001

$ws.Range("H32").Value = 188134.9
$ws.Range("I32").Value = 270439.22
$ws.Range("K32").Value = 270439.22
$ws.Range("M32").Value = -270152.22

$ws.Range("H45").Value = 1825.0834
$ws.Range("I45").Value = 1858.6666
$ws.Range("J45").Value = 1791.5
$ws.Range("K45").Value = 1858.6666
$ws.Range("L45").Value = 1791.5
$ws.Range("M45").Value = -1481.6666
$ws.Range("N45").Value = -2545.5

$ws.Range("H61").Value = 1758141.9
$ws.Range("I61").Value = 3728.46
$ws.Range("K61").Value = 3728.46
$ws.Range("M61").Value = -3516.46

$ws.Range("H132").Value = 628607.8
$ws.Range("I132").Value = 717794.8
$ws.Range("J132").Value = 4298.8
$ws.Range("K132").Value = 2153384.4
$ws.Range("L132").Value = 12896.4
$ws.Range("M132").Value = -2150854.4
$ws.Range("N132").Value = -17956.4

$ws.Range("H134").Value = 52741.7
$ws.Range("J134").Value = 52741.7
$ws.Range("L134").Value = 52741.7
$ws.Range("N134").Value = -62881.7

$ws.Range("H136").Value = 1758141.9
$ws.Range("I136").Value = 3728.46
$ws.Range("K136").Value = 11185.38
$ws.Range("M136").Value = -8635.380000000001

$ws.Range("H138").Value = 104628
$ws.Range("J138").Value = 104628
$ws.Range("L138").Value = 104628
$ws.Range("N138").Value = -114908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 257648.9
$ws.Range("I94").Value = 8880.471
$ws.Range("J94").Value = 1667336.6
$ws.Range("K94").Value = 8880.471
$ws.Range("L94").Value = 1667336.6
$ws.Range("M94").Value = -8429.471
$ws.Range("N94").Value = -1668238.6

$ws.Range("H132").Value = 60000
$ws.Range("J132").Value = 60000
$ws.Range("L132").Value = 60000
$ws.Range("N132").Value = -70120

$ws.Range("H134").Value = 3093540
$ws.Range("I134").Value = 5202.024
$ws.Range("J134").Value = 13902723
$ws.Range("K134").Value = 15606.072
$ws.Range("L134").Value = 41708169
$ws.Range("M134").Value = -13071.072
$ws.Range("N134").Value = -41713239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 133053.58
$ws.Range("I16").Value = 116556
$ws.Range("J16").Value = 168798.33
$ws.Range("K16").Value = 116556
$ws.Range("L16").Value = 168798.33
$ws.Range("M16").Value = -116269
$ws.Range("N16").Value = -169372.33

$ws.Range("H31").Value = 1465213.5
$ws.Range("I31").Value = 2225041
$ws.Range("J31").Value = 4006.7693
$ws.Range("K31").Value = 2225041
$ws.Range("L31").Value = 4006.7693
$ws.Range("M31").Value = -2224746
$ws.Range("N31").Value = -4596.7693

$ws.Range("H34").Value = 1465213.5
$ws.Range("I34").Value = 2225041
$ws.Range("J34").Value = 4006.7693
$ws.Range("K34").Value = 2225041
$ws.Range("L34").Value = 4006.7693
$ws.Range("M34").Value = -2224839
$ws.Range("N34").Value = -4410.7693

$ws.Range("H58").Value = 2535334
$ws.Range("I58").Value = 8027.7617
$ws.Range("J58").Value = 6958119.5
$ws.Range("K58").Value = 8027.7617
$ws.Range("L58").Value = 6958119.5
$ws.Range("M58").Value = -7824.7617
$ws.Range("N58").Value = -6958525.5

$ws.Range("H62").Value = 4471.846
$ws.Range("I62").Value = 4946.5
$ws.Range("K62").Value = 4946.5
$ws.Range("M62").Value = -4322.5

$ws.Range("H65").Value = 4471.846
$ws.Range("I65").Value = 4946.5
$ws.Range("K65").Value = 24732.5
$ws.Range("M65").Value = -21612.5

$ws.Range("H107").Value = 355.6111
$ws.Range("J107").Value = 343.5
$ws.Range("L107").Value = 343.5
$ws.Range("N107").Value = -4183.5

$ws.Range("H113").Value = 133053.58
$ws.Range("I113").Value = 116556
$ws.Range("J113").Value = 168798.33
$ws.Range("K113").Value = 116556
$ws.Range("L113").Value = 168798.33
$ws.Range("M113").Value = -114386
$ws.Range("N113").Value = -173138.33

$ws.Range("H122").Value = 7897.963
$ws.Range("I122").Value = 1855.5769
$ws.Range("K122").Value = 5566.7307
$ws.Range("M122").Value = -3116.7307

$ws.Range("H136").Value = 2535334
$ws.Range("I136").Value = 8027.7617
$ws.Range("J136").Value = 6958119.5
$ws.Range("K136").Value = 24083.2851
$ws.Range("L136").Value = 20874358.5
$ws.Range("M136").Value = -21533.2851
$ws.Range("N136").Value = -20879458.5

$ws.Range("H141").Value = 208885.89
$ws.Range("J141").Value = 221696.06
$ws.Range("L141").Value = 221696.06
$ws.Range("N141").Value = -232056.06

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3171.2856
$ws.Range("J39").Value = 4199.8
$ws.Range("L39").Value = 12599.4
$ws.Range("N39").Value = -13187.4

$ws.Range("H122").Value = 807531.6
$ws.Range("J122").Value = 1320.625
$ws.Range("L122").Value = 11885.625
$ws.Range("N122").Value = -16785.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 14999.667
$ws.Range("J26").Value = 14999.667
$ws.Range("L26").Value = 14999.667
$ws.Range("N26").Value = -15559.667

$ws.Range("H50").Value = 14999.667
$ws.Range("J50").Value = 14999.667
$ws.Range("L50").Value = 14999.667
$ws.Range("N50").Value = -15995.667

$ws.Range("H132").Value = 7909.609
$ws.Range("I132").Value = 6746.6035
$ws.Range("J132").Value = 14041.818
$ws.Range("K132").Value = 20239.8105
$ws.Range("L132").Value = 42125.454
$ws.Range("M132").Value = -17709.8105
$ws.Range("N132").Value = -47185.454

$ws.Range("H134").Value = 47847.1
$ws.Range("J134").Value = 47847.1
$ws.Range("L134").Value = 143541.3
$ws.Range("N134").Value = -148611.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5726.7856
$ws.Range("I7").Value = 6084.1816
$ws.Range("J7").Value = 4416.3335
$ws.Range("K7").Value = 6084.1816
$ws.Range("L7").Value = 4416.3335
$ws.Range("M7").Value = -5972.1816
$ws.Range("N7").Value = -4640.3335

$ws.Range("H43").Value = 1795718.2
$ws.Range("J43").Value = 1795718.2
$ws.Range("L43").Value = 1795718.2
$ws.Range("N43").Value = -1796104.2

$ws.Range("H99").Value = 30259
$ws.Range("I99").Value = 30259
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 30259
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -27264
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 5726.7856
$ws.Range("I126").Value = 6084.1816
$ws.Range("J126").Value = 4416.3335
$ws.Range("K126").Value = 18252.5448
$ws.Range("L126").Value = 13249.0005
$ws.Range("M126").Value = -15782.5448
$ws.Range("N126").Value = -18189.0005

$ws.Range("H132").Value = 1859867.6
$ws.Range("I132").Value = 3185836.5
$ws.Range("J132").Value = 3511
$ws.Range("K132").Value = 9557509.5
$ws.Range("L132").Value = 10533
$ws.Range("M132").Value = -9554979.5
$ws.Range("N132").Value = -15593

$ws.Range("H136").Value = 7582242
$ws.Range("I136").Value = 8336350.5
$ws.Range("J136").Value = 6953818
$ws.Range("K136").Value = 25009051.5
$ws.Range("L136").Value = 20861454
$ws.Range("M136").Value = -25006501.5
$ws.Range("N136").Value = -20866554

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 24501.5
$ws.Range("J103").Value = 24501.5
$ws.Range("L103").Value = 24501.5
$ws.Range("N103").Value = -26845.5

$ws.Range("H132").Value = 3705183
$ws.Range("I132").Value = 4274821.5
$ws.Range("J132").Value = 2533.3333
$ws.Range("K132").Value = 12824464.5
$ws.Range("L132").Value = 7599.999899999999
$ws.Range("M132").Value = -12821934.5
$ws.Range("N132").Value = -12659.9999

$ws.Range("H136").Value = 12174707
$ws.Range("I136").Value = 2718289.5
$ws.Range("J136").Value = 50000376
$ws.Range("K136").Value = 8154868.5
$ws.Range("L136").Value = 150001128
$ws.Range("M136").Value = -8152318.5
$ws.Range("N136").Value = -150006228

$ws.Range("H137").Value = 94080.8
$ws.Range("J137").Value = 94080.8
$ws.Range("L137").Value = 94080.8
$ws.Range("N137").Value = -104280.8
